$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-affecting data: rows 2-4 get new values, rows 5-7 are new.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.

$data = @(
    @("ECs","Il1rapl1","Ptprf","ECs",1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,3,1,0.21127,0.63381,0.02468048274231428,0.02468048274231428,0.005096114093333333,0.04586502684,0.005502039270015466,0.005502039270015467),
    @("ECs","Il1rapl1","Ptprf","FAPs",1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,3,1,3.467027333333334,10.401082,0.405016842275123,0.405016842275123,0.08362932198311111,0.752663897848,0.09029072058606052,0.09029072058606052),
    @("ECs","Il1rapl1","Ptprf","sCs",1,0.3333333333333333,0.02412133333333333,0.072364,0.2229307800605663,0.2229307800605663,3,1,4.881908,14.645724,0.5703026749825627,0.5703026749825627,0.1177581301706667,1.059823171536,0.1271380202044903,0.1271380202044903),
    @("sCs","Il1rapl1","Ptprf","ECs",3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,3,1,0.21127,0.63381,0.02468048274231428,0.02468048274231428,0.01776351117666666,0.15987160059,0.01917844347229881,0.01917844347229881),
    @("sCs","Il1rapl1","Ptprf","FAPs",3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,3,1,3.467027333333334,10.401082,0.405016842275123,0.405016842275123,0.2915065025108889,2.623558522598,0.3147261216890626,0.3147261216890626),
    @("sCs","Il1rapl1","Ptprf","sCs",3,1,0.08407966666666666,0.252239,0.7770692199394338,0.7770692199394338,3,1,4.881908,14.645724,0.5703026749825627,0.5703026749825627,0.4104691973373333,3.694222776036,0.4431646547780724,0.4431646547780724)
)

$rowIdx = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $val
        $colIdx++
    }
    $rowIdx++
}
